# Applies the cryptos-list price/volume refresh described in the commit diff.
# Numeric-looking text values (e.g. "217.29") are written with a leading
# apostrophe so Excel keeps them as text (matching the source t="inlineStr"
# cells) instead of silently converting them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.257.15"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.644.59"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'217.29"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'0.517"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "'20.04"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "'0.0850"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.876.26"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "1.639.04"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'0.545"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "'67.05"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "27.250.77"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'220.10"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'6.99"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("D22").Value = "'2.54"
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("D23").Value = "'4.42"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'148.72"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "'15.74"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.305.87"
$ws.Range("E34").Value = "  +3.76%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.58"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("D39").Value = "'0.859"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").Value = "'0.812"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  +6.10%  "
$ws.Range("D43").Value = "'5.32"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").Value = "1.785.32"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'62.05"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'92.05"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'7.70"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "'0.0967"
$ws.Range("E51").Value = "  +0.24%  "
